$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 13: without-opt M30 result, taken from H1 ---
$ws.Range("A13").Value = "without opt, taken from H1"
$ws.Range("B13").Value = "eur/usd"
$ws.Range("C13").Value = "M30"
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = 90
$ws.Range("G13").Value = 35
$ws.Range("Q13").Value = "profitable on M30 based on H1 result"
$ws.Range("O13").Value = "+"

# --- New row 14: without-opt M15 result, taken from H2 ---
$ws.Range("O14").Value = "-"
$ws.Range("Q14").Value = "lossy on M15 based on H1 result"
$ws.Range("A14").Value = "without opt, taken from H2"
$ws.Range("B14").Value = "eur/usd"
$ws.Range("C14").Value = "M15"
$ws.Range("D14").Value = 12
$ws.Range("E14").Value = 90
$ws.Range("G14").Value = 35

# --- New row 15: opt len and correl on M15 ---
$ws.Range("A15").Value = "opt len and correl"
$ws.Range("Q15").Value = "opt on M15 seems to be unstable"
$ws.Range("B15").Value = "eur/usd"
$ws.Range("C15").Value = "M15"
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = 90
$ws.Range("G15").Value = 35

# --- Rename header M1: "processed bars" -> "to/processed bars" ---
$ws.Range("M1").Value = "to/processed bars"

# --- M15 note & remaining row-15 numbers ---
$ws.Range("M15").Value = "10/2015 - 10000"
$ws.Range("N15").Value = 80
$ws.Range("O15").Value = 670
$ws.Range("P15").Value = 1.37

# L15 needs the same date format/style as the other "from" date cells (L2/L8)
$ws.Range("L8").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("L15").Value = 42005

# --- Row 8: add M8 ("4/2015", text-formatted), clear old Q8 note ---
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "4/2015"
$ws.Range("Q8").ClearContents()

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection state as left by the author after the edit ---
$ws.Range("N16").Select()
